# Update the cash-flow report from "Aug 2022" to "Sep 2022" figures.
#
# The worksheet stores its numbers (even the ones with thousands
# separators) as literal text in shared strings rather than as formatted
# numeric cells, so a plain `.Value = "141,680,583"` assignment would get
# auto-coerced by Excel into a real number (losing the text/comma
# formatting, and picking up an unwanted NumberFormat/style). To keep the
# cells as plain text - matching how the workbook already stores them -
# we stage the literal through a helper cell as a `="..."` formula,
# copy it, and paste-special as values only. That converts the formula
# result into a literal text cell without ever touching NumberFormat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z100")

function Set-TextValue {
    param($rangeAddress, $text)
    $helper.Formula = "=""" + $text + """"
    $helper.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)
}

# Column headers: "Aug 2022" / "SD Aug 2022" -> "Sep 2022" / "SD Sep 2022"
Set-TextValue "D3" "Sep 2022"
Set-TextValue "E3" "SD Sep 2022"

# 1. Pengembalian Pinjaman Mitra Binaan (Sep 2022 column)
Set-TextValue "D6" "141,680,583"

# 1. Penyaluran Pinjaman Kemitraan (Sep 2022 column) now empty/zero
$ws.Range("D12").Value = 0

# KAS NETTO DITERIMA(DIGUNAKAN) UNTUK AKTIVITAS OPERASI (Sep 2022 column) now zero
$ws.Range("D20").Value = 0

# KENAIKAN (PENURUNAN) NETTO DALAM KAS/SETARA KAS
Set-TextValue "D25" "11,867,642"
Set-TextValue "E25" "581,867,642"

# KAS DAN SETARA KAS PADA AWAL TAHUN
Set-TextValue "D26" "756,829,644"
Set-TextValue "E26" "1,458,228,442"

# KAS DAN SETARA KAS PADA AKHIR TAHUN
Set-TextValue "D27" "768,697,286"
Set-TextValue "E27" "2,040,096,084"

# Clean up the helper cell
$helper.Clear()
